$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style (bold, centered, bordered header look) from the existing
# "panel" header cell (E1) onto the new "time_taken" header cell (F1),
# then set its text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Data cells F2:F9 with the recorded timestamps (plain, unstyled like the
# rest of the data rows)
$timestamps = @(
    "2021-10-05 13:39:45.344115",
    "2021-10-05 13:39:45.344125",
    "2021-10-05 13:39:45.344129",
    "2021-10-05 13:39:45.344132",
    "2021-10-05 13:39:45.344135",
    "2021-10-05 13:39:45.344137",
    "2021-10-05 13:39:45.344140",
    "2021-10-05 13:39:45.344143"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
